# Workbook: departure-procedure distance analysis (Radar / Wake-turbulence).
# Fills in the previously-blank "ProcDesp" (departure procedure) values for
# the rows that had a placeholder "-" in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$updates = @{
    11  = "LOTOS3R"
    12  = "NATPI2R"
    13  = "AGENA4R"
    14  = "DIPES1R"
    16  = "LAPRA4R"
    17  = "OLOXO1R"
    18  = "LOTOS3R"
    19  = "OLOXO1R"
    20  = "AGENA4R"
    27  = "LOTOS6Q"
    86  = "LOBAR7Q"
    94  = "LOBAR7Q"
    97  = "GRAUS6Q"
    102 = "SENIA6Q"
    103 = "NATPI3Q"
    106 = "GRAUS6Q"
    144 = "AGENA6Q"
    212 = "DALIN5Q"
    220 = "DALIN5Q"
    242 = "AGENA6Q"
    252 = "LOTOS6Q"
    255 = "SENIA6Q"
    258 = "SENIA6Q"
    283 = "DALIN5Q"
    292 = "DUNES6Q"
    313 = "LOTOS6Q"
    319 = "LOTOS6Q"
    320 = "LOTOS6Q"
    321 = "LOTOS6Q"
    335 = "LOBAR7Q"
    340 = "DALIN5Q"
    341 = "LOTOS6Q"
    423 = "LAPRA6Q"
    435 = "DIPES2Q"
    438 = "DIPES2Q"
    444 = "DIPES2Q"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}

# Reflect the view state the author left the sheet in: scrolled down,
# zoomed to 80%, with F448 as the active selection.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 163
$win.ScrollColumn = 2
$win.Zoom = 80
$ws.Range("F448").Select()
